$wb = $excel.ActiveWorkbook

# --- Sheet1: get_fwd_release_delays ---
$ws1 = $wb.Worksheets.Item("get_fwd_release_delays")
$ws1.Range("A1").Value = 1
$ws1.Range("B1").Value = 1
$ws1.Range("A2").Value = 2
$ws1.Range("B2").Value = 2
$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = 2
[void]$ws1.Range("B1").Select()

# --- Sheet2: get_fwd_proc_compute_node ---
$ws2 = $wb.Worksheets.Item("get_fwd_proc_compute_node")
$ws2.Range("A1").Value = 2
$ws2.Range("A2").Value = 4
$ws2.Range("B1").ClearContents()
[void]$ws2.Range("A3").Select()

# --- Sheet3: get_fwd_end_local ---
$ws3 = $wb.Worksheets.Item("get_fwd_end_local")
$ws3.Range("A1").Value = 1
$ws3.Range("A2").Value = 5
$ws3.Range("A3").Value = 4
$ws3.Range("B1").ClearContents()
$ws3.Range("C1").ClearContents()
[void]$ws3.Range("C5").Select()

# --- Sheet5: get_memory_characteristics ---
$ws5 = $wb.Worksheets.Item("get_memory_characteristics")
$ws5.Range("A1").Value = 3
$ws5.Range("A2").Value = 10
$ws5.Range("B1").ClearContents()
[void]$ws5.Range("A2").Select()
[void]$ws5.Activate()
